# Excel_Challenge_608 - Extract Zip and Country
# Commit: "I think I finally understand what is going on"
#
# Summary of edits applied to the "EDA" worksheet:
#  - Insert one new (blank, Heading-1-styled) row above the old row 57,
#    pushing the "single-cell REGEXEXTRACT" demo block (old rows 57-66)
#    down to rows 58-67.
#  - Give column F an explicit width.
#  - Add a new "Alternate Formula Form" section below that block:
#      * two "Comment" styled notes
#      * a "Heading 1" styled caption ("Alternate Formula Form")
#      * a new HSTACK/REDUCE/LAMBDA dynamic-array formula (F74:G83)
#        that recomputes the Zip/Country columns in one shot, plus an
#        ANCHORARRAY comparison against B3:C12 (I74:J83)
#  - Update the view (drop the old scrolled/selected position, select I60)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDA")
$ws.Activate()

# ---------------------------------------------------------------------
# 1) Insert a new blank row above the old row 57 (shifts old rows
#    57-66 down to 58-67). Excel's default "insert" behaviour copies
#    the formatting of the row above (row 56, which uses the
#    "Heading 1" cell style / s="8"), which is exactly what the target
#    file shows for the newly inserted F57.
# ---------------------------------------------------------------------
$ws.Rows.Item(57).Insert()

# ---------------------------------------------------------------------
# 2) Column F width
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 13.43

# ---------------------------------------------------------------------
# 3) New "Alternate Formula Form" block
# ---------------------------------------------------------------------

# Two "Comment" styled notes (rows 69-70)
$ws.Cells.Item(69, 6).Style = "Comment"
$ws.Cells.Item(69, 6).Value = "Using a 0 or more (*) operator results in two matches and the function showed the last one, which was empty."

$ws.Cells.Item(70, 6).Style = "Comment"
$ws.Cells.Item(70, 6).Value = "Why was there an an empty match? It was an option that was TRUE. Using the + operator removed the empty option."

# "Heading 1" styled caption (row 72)
$ws.Cells.Item(72, 6).Style = "Heading 1"
$ws.Cells.Item(72, 6).Value = "Alternate Formula Form"

# New dynamic-array formula block (rows 74-83)
$ws.Range("F74:G83").FormulaArray = '=_xlfn.HSTACK(0+_xlfn.DROP(_xlfn.REDUCE("",A3:A12, _xlfn.LAMBDA(_xlpm.a,_xlpm.v,_xlfn.VSTACK(_xlpm.a,_xlfn.TAKE(_xlfn.REGEXEXTRACT(_xlpm.v,"\d{5,6}",1),,-1)))),1),TRIM(_xlfn.REGEXEXTRACT(A3:A12,"(United |New )?[A-Za-z]+$",1)))'

$ws.Range("I74:J83").FormulaArray = '=_xlfn.ANCHORARRAY(F74)=B3:C12'

# ---------------------------------------------------------------------
# 4) View: clear the old scrolled position, select I60
# ---------------------------------------------------------------------
$ws.Range("I60").Select()
